$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (Excel shifts rows 8..16 down to 9..17,
# preserving per-row formatting from the row pushed down).
$ws.Rows(8).Insert()

# The literal serial numbers in column A (1,2,3,...) need to be bumped by
# one for every data row that moved down, since the new "建筑公司" entry
# now occupies position 7 in that sequence.
for ($r = 17; $r -ge 9; $r--) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value + 1
}

# Fill in the new "建筑公司" (Construction Company) project row, mirroring
# the layout used by the other project-style row (核能发电站, row 4):
# only A/B/H/I/J/L are populated.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "建筑公司"
$ws.Range("H8").Value = 4000
$ws.Range("I8").Value = 3000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = "每当玩家升级地产时，获得500元。当任意玩家到建筑公司时，可将一处地产升一级（需支付升级费用）。"

# Match L4's style (wrap text, vertically centered) and the row's authored height.
$ws.Range("L8").Style = $ws.Range("L4").Style
$ws.Rows(8).RowHeight = 42.75

# Update the active selection to match the post-edit state.
$ws.Range("F8").Select()
